# Update leve profit calculations across all profession sheets
# (scheduled runner refresh of currentAveragePrice / LevePrice / LeveProfit columns)
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2236.3635
$ws.Range("I19").Value = 1049.3334
$ws.Range("J19").Value = 2681.5
$ws.Range("K19").Value = 1049.3334
$ws.Range("L19").Value = 2681.5
$ws.Range("M19").Value = -874.3334
$ws.Range("N19").Value = -3031.5
$ws.Range("H33").Value = 109
$ws.Range("I33").Value = 115.8
$ws.Range("J33").Value = 92
$ws.Range("K33").Value = 115.8
$ws.Range("L33").Value = 92
$ws.Range("M33").Value = 113.2
$ws.Range("N33").Value = -550
$ws.Range("H48").Value = 2159
$ws.Range("J48").Value = 2159
$ws.Range("L48").Value = 6477
$ws.Range("N48").Value = -7061
$ws.Range("H56").Value = 2159
$ws.Range("J56").Value = 2159
$ws.Range("L56").Value = 6477
$ws.Range("N56").Value = -7545
$ws.Range("H98").Value = 1635.6522
$ws.Range("I98").Value = 1669.4736
$ws.Range("J98").Value = 1475
$ws.Range("K98").Value = 1669.4736
$ws.Range("L98").Value = 1475
$ws.Range("M98").Value = -171.4736
$ws.Range("N98").Value = -4471
$ws.Range("H116").Value = 10422.556
$ws.Range("I116").Value = 21599
$ws.Range("J116").Value = 6123.923
$ws.Range("K116").Value = 21599
$ws.Range("L116").Value = 6123.923
$ws.Range("M116").Value = -18157
$ws.Range("N116").Value = -13007.923
$ws.Range("H122").Value = 1635.6522
$ws.Range("I122").Value = 1669.4736
$ws.Range("J122").Value = 1475
$ws.Range("K122").Value = 5008.4208
$ws.Range("L122").Value = 4425
$ws.Range("M122").Value = -2558.4208
$ws.Range("N122").Value = -9325
$ws.Range("H138").Value = 2884.8438
$ws.Range("J138").Value = 3000
$ws.Range("L138").Value = 9000
$ws.Range("N138").Value = -19280

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 926576.3
$ws.Range("J2").Value = 794.5
$ws.Range("L2").Value = 794.5
$ws.Range("N2").Value = -1020.5
$ws.Range("H24").Value = 83038.8
$ws.Range("J24").Value = 83038.8
$ws.Range("L24").Value = 83038.8
$ws.Range("N24").Value = -83786.8
$ws.Range("H32").Value = 3533.4807
$ws.Range("J32").Value = 4640.3
$ws.Range("L32").Value = 4640.3
$ws.Range("N32").Value = -5214.3
$ws.Range("H45").Value = 1686
$ws.Range("I45").Value = 1544.5
$ws.Range("K45").Value = 1544.5
$ws.Range("M45").Value = -1167.5
$ws.Range("H61").Value = 2750
$ws.Range("I61").Value = 1893.8077
$ws.Range("K61").Value = 1893.8077
$ws.Range("M61").Value = -1681.8077
$ws.Range("H74").Value = 1129.2703
$ws.Range("J74").Value = 1631.909
$ws.Range("L74").Value = 1631.909
$ws.Range("N74").Value = -3379.909
$ws.Range("H76").Value = 23924.572
$ws.Range("J76").Value = 23924.572
$ws.Range("L76").Value = 23924.572
$ws.Range("N76").Value = -24600.572
$ws.Range("H77").Value = 1129.2703
$ws.Range("J77").Value = 1631.909
$ws.Range("L77").Value = 8159.545
$ws.Range("N77").Value = -16895.545
$ws.Range("H79").Value = 23924.572
$ws.Range("J79").Value = 23924.572
$ws.Range("L79").Value = 23924.572
$ws.Range("N79").Value = -26264.572
$ws.Range("H82").Value = 72220.336
$ws.Range("J82").Value = 73331.2
$ws.Range("L82").Value = 73331.2
$ws.Range("N82").Value = -74053.2
$ws.Range("H85").Value = 72220.336
$ws.Range("J85").Value = 73331.2
$ws.Range("L85").Value = 73331.2
$ws.Range("N85").Value = -75827.2
$ws.Range("H100").Value = 83038.8
$ws.Range("J100").Value = 83038.8
$ws.Range("L100").Value = 83038.8
$ws.Range("N100").Value = -85202.8
$ws.Range("H112").Value = 45000
$ws.Range("J112").Value = 45000
$ws.Range("L112").Value = 45000
$ws.Range("N112").Value = -47954
$ws.Range("H116").Value = 926576.3
$ws.Range("J116").Value = 794.5
$ws.Range("L116").Value = 794.5
$ws.Range("N116").Value = -5382.5
$ws.Range("H122").Value = 2303.4285
$ws.Range("I122").Value = 2103.6924
$ws.Range("K122").Value = 6311.0772
$ws.Range("M122").Value = -3861.0772
$ws.Range("H125").Value = 100000
$ws.Range("J125").Value = 100000
$ws.Range("L125").Value = 100000
$ws.Range("N125").Value = -109840
$ws.Range("H132").Value = 1848.0278
$ws.Range("I132").Value = 1316.125
$ws.Range("K132").Value = 3948.375
$ws.Range("M132").Value = -1418.375
$ws.Range("H136").Value = 2750
$ws.Range("I136").Value = 1893.8077
$ws.Range("K136").Value = 5681.4231
$ws.Range("M136").Value = -3131.4231

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 926576.3
$ws.Range("J3").Value = 794.5
$ws.Range("L3").Value = 794.5
$ws.Range("N3").Value = -1022.5
$ws.Range("H99").Value = 2399.75
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 2799.5
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 2799.5
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -5795.5
$ws.Range("H134").Value = 5363.4873
$ws.Range("J134").Value = 2672.5454
$ws.Range("L134").Value = 8017.6362
$ws.Range("N134").Value = -13087.6362

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2008.9131
$ws.Range("I31").Value = 1752.6666
$ws.Range("K31").Value = 1752.6666
$ws.Range("M31").Value = -1457.6666
$ws.Range("H34").Value = 2008.9131
$ws.Range("I34").Value = 1752.6666
$ws.Range("K34").Value = 1752.6666
$ws.Range("M34").Value = -1550.6666
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H105").Value = 2255.125
$ws.Range("I105").Value = 1836.5
$ws.Range("J105").Value = 3511
$ws.Range("K105").Value = 1836.5
$ws.Range("L105").Value = 3511
$ws.Range("M105").Value = -89.5
$ws.Range("N105").Value = -7005
$ws.Range("H134").Value = 1678.0714
$ws.Range("I134").Value = 1468.9
$ws.Range("J134").Value = 2201
$ws.Range("K134").Value = 4406.700000000001
$ws.Range("L134").Value = 6603
$ws.Range("M134").Value = -1871.700000000001
$ws.Range("N134").Value = -11673

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7013.8125
$ws.Range("I56").Value = 7013.8125
$ws.Range("K56").Value = 7013.8125
$ws.Range("M56").Value = -6483.8125

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 32981.816
$ws.Range("J135").Value = 32981.816
$ws.Range("L135").Value = 32981.816
$ws.Range("N135").Value = -43121.816

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 11656.333
$ws.Range("J104").Value = 11656.333
$ws.Range("L104").Value = 11656.333
$ws.Range("N104").Value = -18644.333
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 1578.5538
$ws.Range("I132").Value = 1396.4773
$ws.Range("K132").Value = 4189.4319
$ws.Range("M132").Value = -1659.4319

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H119").Value = 26088.2
$ws.Range("J119").Value = 26088.2
$ws.Range("L119").Value = 26088.2
$ws.Range("N119").Value = -35764.2
$ws.Range("H136").Value = 12921498
$ws.Range("I136").Value = 13890310
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 41670930
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -41668380
$ws.Range("N136").Value = -17100
